# Applies the cryptos-list price/volume refresh described in the commit
# "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.818.43"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.44%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.634.88"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.39%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.26%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.29%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5020"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.45%  "

# Row 7
$ws.Range("E7").Value = "  -0.26%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2568"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.59%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06407"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.43%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.59"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.78%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07670"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.63%  "

# Row 12
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.636.39"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.16%  "

# Row 13
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.230"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.44%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.858.65"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.44%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5452"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.62%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅7914"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.84%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.48"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.16%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.851.47"
$ws.Range("D18").Style = "Normal"

# Row 19
$ws.Range("E19").Value = "  -0.18%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "202.87"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.37%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.300"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.80%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.930"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.03%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.969"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.20%  "

# Row 24
$ws.Range("E24").Value = "  -0.18%  "

# Row 25
$ws.Range("E25").Value = "  +10.26%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "140.97"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.67%  "

# Row 27
$ws.Range("E27").Value = "  -1.70%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.68"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.75%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.689"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.00%  "

# Row 30
$ws.Range("E30").Value = "  -1.20%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.04977"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.77%  "

# Row 32
$ws.Range("E32").Value = "  -2.73%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.182"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.14%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.530"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.37%  "

# Row 35
$ws.Range("E35").Value = "  -0.70%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.171.21"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.26%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.624"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.93%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.8912"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.09%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5566"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.99%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01555"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.27%  "

# Row 41
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.003"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.21%  "

# Row 42
$ws.Range("B42").Value = "mCoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.546"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.55%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.629"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.71%  "

# Row 44
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.26"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.14%  "

# Row 45
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8016"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.14%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.771.50"
$ws.Range("D46").Style = "Normal"

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₈115"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.47%  "

# Row 48
$ws.Range("E48").Value = "  -0.58%  "

# Row 49
$ws.Range("E49").Value = "  +0.16%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "54.74"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.08%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05031"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.47%  "

